$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name and workbook title to reflect the new "through" date
$ws.Name = "Through 2021-11-20"

# Update the row label for November to reflect the new "through" date
$ws.Range("A12").Value = "November (through 11-20)"

# Update November row (row 12) values
$ws.Range("B12").Value = 21
$ws.Range("C12").Value = 46
$ws.Range("D12").Value = 82
$ws.Range("E12").Value = 42
$ws.Range("F12").Value = 30
$ws.Range("G12").Value = 127
$ws.Range("H12").Value = 137

# Update Total row (row 13) values
$ws.Range("B13").Value = 279
$ws.Range("C13").Value = 532
$ws.Range("D13").Value = 792
$ws.Range("E13").Value = 657
$ws.Range("F13").Value = 512
$ws.Range("G13").Value = 1184
$ws.Range("H13").Value = 1579
